# saitenuni.xlsx update
# 1) Bump the "login count" / "entry count" tallies for several universities
#    on Sheet1 (and the grand-total row), matching a refreshed data pull.
# 2) Add two new blank worksheets (Sheet2, Sheet3) after Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated counts (column C = ログイン者(人), column D = 入力者(人)) ---
$ws.Range("C16").Value = 75
$ws.Range("D19").Value = 37
$ws.Range("C24").Value = 68
$ws.Range("D24").Value = 63
$ws.Range("D29").Value = 54
$ws.Range("C31").Value = 62
$ws.Range("D31").Value = 48
$ws.Range("C33").Value = 50
$ws.Range("D33").Value = 36
$ws.Range("C35").Value = 57
$ws.Range("D35").Value = 49
$ws.Range("C37").Value = 59
$ws.Range("D37").Value = 42
$ws.Range("C40").Value = 79
$ws.Range("D40").Value = 69
$ws.Range("D43").Value = 53
$ws.Range("C46").Value = 68
$ws.Range("C49").Value = 61
$ws.Range("D49").Value = 53
$ws.Range("C52").Value = 50
$ws.Range("C64").Value = 30
$ws.Range("D64").Value = 23
$ws.Range("C68").Value = 56
$ws.Range("D68").Value = 44
$ws.Range("D74").Value = 94
$ws.Range("C77").Value = 132
$ws.Range("D77").Value = 131
$ws.Range("C81").Value = 74
$ws.Range("C84").Value = 188
$ws.Range("C92").Value = 243

# --- Grand total row ---
$ws.Range("C93").Value = 5428
$ws.Range("D93").Value = 4358

# --- Add two new blank worksheets after Sheet1 (Sheet2, Sheet3) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet2 = $wb.Worksheets.Add($null, $lastSheet)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $lastSheet)

# Keep Sheet1 as the active / selected tab, as in the source workbook.
$wb.Worksheets.Item("Sheet1").Activate()
